$d = $word.ActiveDocument

# Locate the "Requisitos" paragraph that mentions LOB1009. Immediately
# after it there is an empty paragraph, a "Ver no Jupiter..." paragraph,
# and a "(c) 2020 ..." footer paragraph, all of which were removed from
# the page (the trailing blank paragraph before the page-break paragraph
# is kept as-is).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*LOB1009*") {
        $target = $i
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the LOB1009 paragraph"
}

$startPara = $d.Paragraphs.Item($target + 1)
$endPara = $d.Paragraphs.Item($target + 3)
$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
$rng.Delete()
